$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Transmitance values in column B (rows 3-18) from 1 to 100
$ws.Range("B3:B18").Value = 100

# Update the active cell selection to B18, matching the saved selection state
$ws.Range("B18").Select()
